$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 122: Rio de Janeiro night 1 - guitar
$ws.Range("E122").Value = "yellow"
$ws.Range("F122").Value = "guitar"
$ws.Range("G122").Value = "Stay Beautiful"

# Row 123: Rio de Janeiro night 1 - piano
$ws.Range("E123").Value = "yellow"
$ws.Range("F123").Value = "piano"
$ws.Range("G123").Value = "Suburban Legends (Taylor's Version) [From The Vault]"

# Row 124: Rio de Janeiro night 2 - guitar (date shifts 45248 -> 45249)
$ws.Range("B124").Value = 45249
$ws.Range("E124").Value = "blue"
$ws.Range("F124").Value = "guitar"
$ws.Range("G124").Value = "Dancing With Our Hands Tied"

# Row 125: Rio de Janeiro night 2 - piano (date shifts 45248 -> 45249)
$ws.Range("B125").Value = 45249
$ws.Range("E125").Value = "blue"
$ws.Range("F125").Value = "piano"
$ws.Range("G125").Value = "Bigger Than The Whole Sky"

# Row 126: Rio de Janeiro night 3 - guitar (date shifts 45249 -> 45250)
$ws.Range("B126").Value = 45250
$ws.Range("E126").Value = "green"
$ws.Range("F126").Value = "guitar"
$ws.Range("G126").Value = "ME!"

# Row 127: Rio de Janeiro night 3 - piano (date shifts 45249 -> 45250)
$ws.Range("B127").Value = 45250
$ws.Range("E127").Value = "green"
$ws.Range("F127").Value = "piano"
$ws.Range("G127").Value = "So It Goes..."

# Update the view state to match: topLeftCell C103, selection G123
$ws.Application.ActiveWindow.ScrollRow = 103
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("G123").Select() | Out-Null
